$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.109361886978149
$ws.Range("B1").Value = 2.828197717666626
$ws.Range("C1").Value = 6.910793781280518
$ws.Range("D1").Value = 2.02814245223999
$ws.Range("E1").Value = 1.082121253013611
